# Hillcrest Dairy / COLLEGETOWN order — append two new line items
# (chocolate milk, and an additional 2% milk entry) to the order sheet.
# Quantity / price columns hold text-like values in this workbook (e.g.
# "23", "$16.08"), so force Text formatting on the new numeric-looking
# cells before assigning them — otherwise Excel auto-coerces "2" and
# "$16.20" into real numbers/currency.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5:E6").NumberFormat = "@"

$ws.Range("A5").Value = "chop"
$ws.Range("B5").Value = "Milk - Chocolate (12/16oz)"
$ws.Range("C5").Value = "2"
$ws.Range("D5").Value = "$16.20"
$ws.Range("E5").Value = "$32.40"

$ws.Range("A6").Value = "twop"
$ws.Range("B6").Value = "Milk - 2% (12/16oz)"
$ws.Range("C6").Value = "1"
$ws.Range("D6").Value = "$15.00"
$ws.Range("E6").Value = "$15.00"
